$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '26.947.39'
$ws.Range('E2').Value = '  -3.58%  '

# Row 3
$ws.Range('D3').Value = '1.714.62'
$ws.Range('E3').Value = '  -3.07%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.03%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.01'
$ws.Range('E5').Value = '  -6.39%  '

# Row 6
$ws.Range('E6').Value = '  +0.04%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4784'
$ws.Range('E7').Value = '  +4.40%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3472'
$ws.Range('E8').Value = '  -1.54%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '41.91'
$ws.Range('E9').Value = '  -0.45%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07231'
$ws.Range('E10').Value = '  -2.23%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.040'
$ws.Range('E11').Value = '  -5.12%  '

# Row 12
$ws.Range('E12').Value = '  -0.03%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.82'
$ws.Range('E13').Value = '  -4.35%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.827'
$ws.Range('E14').Value = '  -2.98%  '

# Row 15
$ws.Range('D15').Value = '1.714.34'
$ws.Range('E15').Value = '  -3.27%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.832'
$ws.Range('E16').Value = '  -5.01%  '

# Row 17
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001036'
$ws.Range('E17').Value = '  -2.29%  '

# Row 18
$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '86.16'
$ws.Range('E18').Value = '  -7.01%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06375'
$ws.Range('E19').Value = '  -1.15%  '

# Row 20
$ws.Range('E20').Value = '  +0.01%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.47'
$ws.Range('E21').Value = '  -2.71%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.606'
$ws.Range('E22').Value = '  -2.89%  '

# Row 23
$ws.Range('D23').Value = '27.003.66'
$ws.Range('E23').Value = '  -3.46%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.73'
$ws.Range('E24').Value = '  -4.29%  '

# Row 25
$ws.Range('E25').Value = '  -2.93%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '152.16'
$ws.Range('E26').Value = '  -5.66%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.90'
$ws.Range('E27').Value = '  -1.36%  '

# Row 28
$ws.Range('D28').Value = '1.908.36'

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.074'
$ws.Range('E29').Value = '  -4.17%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '120.87'
$ws.Range('E30').Value = '  -2.61%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.025'
$ws.Range('E31').Value = '  -4.80%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09132'
$ws.Range('E32').Value = '  -1.77%  '

# Row 33
$ws.Range('E33').Value = '  -1.86%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.318'
$ws.Range('E34').Value = '  -5.01%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.467'
$ws.Range('E35').Value = '  +6.36%  '

# Row 36
$ws.Range('E36').Value = '  -4.41%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05853'
$ws.Range('E37').Value = '  -4.46%  '

# Row 38
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2000'
$ws.Range('E38').Value = '  -4.15%  '

# Row 39
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.6052'
$ws.Range('E39').Value = '  -3.35%  '

# Row 40
$ws.Range('B40').Value = 'Aptos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '10.93'
$ws.Range('E40').Value = '  -7.72%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.717'
$ws.Range('E41').Value = '  -4.50%  '

# Row 42
$ws.Range('E42').Value = '  -8.11%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.420'
$ws.Range('E43').Value = '  -5.06%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '12.70'
$ws.Range('E44').Value = '  -4.04%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.561'
$ws.Range('E45').Value = '  -4.70%  '

# Row 46
$ws.Range('E46').Value = '  -4.09%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '118.69'
$ws.Range('E47').Value = '  -3.07%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.829'
$ws.Range('E48').Value = '  -5.58%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.106'
$ws.Range('E49').Value = '  -2.10%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06648'
$ws.Range('E50').Value = '  -2.51%  '

# Row 51
$ws.Range('E51').Value = '  +0.04%  '
